{"js": "// Office.js (Word JavaScript API) edit script.\n// Applies:\n//   1. \"Remaining Budget: 400.0\" -> \"Remaining Budget: 1250.0\"\n//   2. Replace the two \"Description: X, Amount: Y\" paragraphs with a single\n//      paragraph containing a pretty-printed ASCII table (prettytable style).\n\nconst body = context.document.body;\n\n// --- 1. Update the remaining budget line -------------------------------\nconst budgetResults = body.search(\"Remaining Budget: 400.0\", { matchCase: true });\nbudgetResults.load(\"items\");\nawait context.sync();\n\nif (budgetResults.items.length > 0) {\n  budgetResults.items[0].insertText(\"Remaining Budget: 1250.0\", \"Replace\");\n  await context.sync();\n}\n\n// --- 2. Replace the expense detail paragraphs with a pretty table -------\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nparagraphs.items.forEach((p) => p.load(\"text\"));\nawait context.sync();\n\nlet rentParagraph = null;\nlet phoneParagraph = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const p = paragraphs.items[i];\n  const text = p.text || \"\";\n  if (text.indexOf(\"Description: rent\") >= 0) {\n    rentParagraph = p;\n  } else if (text.indexOf(\"Description: phone bill\") >= 0) {\n    phoneParagraph = p;\n  }\n}\n\nif (!rentParagraph) {\n  throw new Error(\"Could not find the 'Description: rent' paragraph\");\n}\n\nconst tableLines = [\n  \"+-------------+--------+\",\n  \"| Description | Amount |\",\n  \"+-------------+--------+\",\n  \"|     rent    | 1200.0 |\",\n  \"|    phone    |  50.0  |\",\n  \"|     car     | 500.0  |\",\n  \"|   grocery   | 400.0  |\",\n  \"|    other    | 600.0  |\",\n  \"+-------------+--------+\",\n];\n// Joining with \"\\v\" (vertical tab) produces a single run containing\n// alternating <w:t>/<w:br/> children, matching the target line-break markup.\nconst tableText = tableLines.join(\"\\v\");\n\n// Remove the second (phone bill) paragraph entirely; the table text takes\n// its place inside the first paragraph.\nif (phoneParagraph) {\n  phoneParagraph.delete();\n}\n\nrentParagraph.clear();\nrentParagraph.insertText(tableText, \"Start\");\n\nawait context.sync();\n", "ps1": "# Word COM interop (PowerShell-style) edit script.\n# Applies:\n#   1. \"Remaining Budget: 400.0\" -> \"Remaining Budget: 1250.0\"\n#   2. Replace the two \"Description: X, Amount: Y\" paragraphs with a single\n#      paragraph containing a pretty-printed ASCII table (prettytable style).\n\n$d = $word.ActiveDocument\n\n# --- 1. Update the remaining budget line --------------------------------\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"Remaining Budget: 400.0\"\n$find.Replacement.Text = \"Remaining Budget: 1250.0\"\n$find.Execute([ref]$null, [ref]$true, [ref]$null, [ref]$null, [ref]$null, [ref]$null, [ref]$null, [ref]$null, [ref]$null, [ref]$null, 2) | Out-Null\n\n# --- 2. Replace the expense detail paragraphs with a pretty table -------\n$rentParagraph = $null\n$phoneParagraph = $null\nforeach ($p in $d.Paragraphs) {\n    $text = $p.Range.Text\n    if ($text -like \"*Description: rent*\") {\n        $rentParagraph = $p\n    } elseif ($text -like \"*Description: phone bill*\") {\n        $phoneParagraph = $p\n    }\n}\n\nif ($null -eq $rentParagraph) {\n    throw \"Could not find the 'Description: rent' paragraph\"\n}\n\n# Remove the second (phone bill) paragraph entirely; the table text takes\n# its place inside the first paragraph.\nif ($null -ne $phoneParagraph) {\n    $phoneParagraph.Range.Delete()\n}\n\n$lineBreak = [char]11\n$tableLines = @(\n    \"+-------------+--------+\",\n    \"| Description | Amount |\",\n    \"+-------------+--------+\",\n    \"|     rent    | 1200.0 |\",\n    \"|    phone    |  50.0  |\",\n    \"|     car     | 500.0  |\",\n    \"|   grocery   | 400.0  |\",\n    \"|    other    | 600.0  |\",\n    \"+-------------+--------+\"\n)\n# Joining with Chr(11) (vertical tab) produces a single run containing\n# alternating <w:t>/<w:br/> children, matching the target line-break markup.\n$tableText = [string]::Join($lineBreak, $tableLines)\n\n$rentParagraph.Range.Text = $tableText\n"}
